$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.81%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.40%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.197"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.35%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06965"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.07%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.437"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.23%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.552"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.04%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.403"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.11%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9003"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-4.12%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1608"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.42%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07731"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'16.54%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07721"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.71%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02931"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.31%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09007"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.31%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001600"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.09%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0006480"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.49%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006536"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.54%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.471"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.37%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.12%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3234"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.83%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1336"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.30%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.038"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.27%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1598"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'4.94%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04528"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.40%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'2.78%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004144"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.63%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001168"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-6.18%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'3.50%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04368"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.14%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006926"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.01%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1244"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.56%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002067"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.71%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01163"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-4.61%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005827"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'4.27%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("E47").Style = "Normal"
